# Fruta / hortaliza, semanal
# Insert a new weekly record as row 28 in the "Ciruela" sheet, pushing
# all subsequent rows (old 28..90) down by one (new 29..91).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28 (shifts existing rows 28-90 down to 29-91)
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record
$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(28, 3).Value = "Ñuble"
$ws.Cells.Item(28, 4).Value = 44998
$ws.Cells.Item(28, 5).Value = 16
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100103
$ws.Cells.Item(28, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(28, 9).Value = 100103002
$ws.Cells.Item(28, 10).Value = "Ciruela"
$ws.Cells.Item(28, 11).Value = "Angeleno"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 80
$ws.Cells.Item(28, 14).Value = 11000
$ws.Cells.Item(28, 15).Value = 12000
$ws.Cells.Item(28, 16).Value = 11500
$ws.Cells.Item(28, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(28, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(28, 19).Value = 639
$ws.Cells.Item(28, 20).Value = 18

Write-Host "Row 28 inserted and populated"
